{"js": "// The document contains a HYPERLINK field whose display text duplicates\n// its target: {HYPERLINK \"http://qrgames.org\"}http://qrgames.org{END}.\n// The author removed the \"http://\" scheme from the *visible* text only,\n// leaving the field code / actual link target untouched, so the reader\n// sees \"qrgames.org\" instead of \"http://qrgames.org\".\n\nconst body = context.document.body;\n\n// Find the visible run of text that shows the link (\"http://qrgames.org\").\n// Word's `search` only matches displayed text, not hidden field-code\n// instructions, so this correctly targets the field's display/result run\n// and not the \" HYPERLINK \"http://qrgames.org\" \" field instruction text.\nconst results = body.search(\"http://qrgames.org\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the displayed link text \"http://qrgames.org\".');\n}\n\n// Replace only the displayed \"http://\" prefix; keep the rest (\"qrgames.org\").\nconst linkDisplay = results.items[0];\nlinkDisplay.insertText(\"qrgames.org\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document contains a HYPERLINK field whose displayed text duplicates\n# its target: {HYPERLINK \"http://qrgames.org\"}http://qrgames.org{END}.\n# The author removed the \"http://\" scheme from the *visible* text only,\n# leaving the field code / actual link target untouched, so the reader\n# sees \"qrgames.org\" instead of \"http://qrgames.org\".\n\n$d = $word.ActiveDocument\n\n# Locate the HYPERLINK field that points at qrgames.org.\n$targetField = $null\nforeach ($f in $d.Fields) {\n    if ($f.Code.Text -match \"HYPERLINK\" -and $f.Code.Text -match \"qrgames\\.org\") {\n        $targetField = $f\n        break\n    }\n}\n\nif ($targetField -eq $null) {\n    throw \"Could not find the qrgames.org HYPERLINK field.\"\n}\n\n# Only touch the field's displayed result text (not its code, so the\n# actual hyperlink target keeps its http:// scheme) and strip the\n# \"http://\" prefix from what is shown to the reader.\n$resultRange = $targetField.Result\n$find = $resultRange.Find\n$find.ClearFormatting()\n$find.Text = \"http://qrgames.org\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"qrgames.org\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, \"qrgames.org\", 2) | Out-Null\n"}
